$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 1016.8182
$ws.Range("I62").Value = 910.625
$ws.Range("J62").Value = 1300
$ws.Range("K62").Value = 910.625
$ws.Range("L62").Value = 1300
$ws.Range("M62").Value = -286.625
$ws.Range("N62").Value = -2548

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 1016.8182
$ws.Range("I65").Value = 910.625
$ws.Range("J65").Value = 1300
$ws.Range("K65").Value = 4553.125
$ws.Range("L65").Value = 6500
$ws.Range("M65").Value = -1433.125
$ws.Range("N65").Value = -12740

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4437.4546
$ws.Range("I86").Value = 4282.4
$ws.Range("J86").Value = 4566.6665
$ws.Range("K86").Value = 4282.4
$ws.Range("L86").Value = 4566.6665
$ws.Range("M86").Value = -3159.4
$ws.Range("N86").Value = -6812.6665

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 4437.4546
$ws.Range("I89").Value = 4282.4
$ws.Range("J89").Value = 4566.6665
$ws.Range("K89").Value = 21412
$ws.Range("L89").Value = 22833.3325
$ws.Range("M89").Value = -15796
$ws.Range("N89").Value = -34065.3325

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 1989.0834
$ws.Range("I111").Value = 2551
$ws.Range("J111").Value = 1202.4
$ws.Range("K111").Value = 7653
$ws.Range("L111").Value = 3607.2
$ws.Range("M111").Value = -4586
$ws.Range("N111").Value = -9741.200000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 11496984
$ws.Range("I138").Value = 25644042
$ws.Range("J138").Value = 2500
$ws.Range("K138").Value = 76932126
$ws.Range("L138").Value = 7500
$ws.Range("M138").Value = -76926986
$ws.Range("N138").Value = -17780

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1708.7
$ws.Range("I45").Value = 1650.5883
$ws.Range("K45").Value = 1650.5883
$ws.Range("M45").Value = -1273.5883

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 5557458.5
$ws.Range("I132").Value = 6758421.5
$ws.Range("J132").Value = 3005.25
$ws.Range("K132").Value = 20275264.5
$ws.Range("L132").Value = 9015.75
$ws.Range("M132").Value = -20272734.5
$ws.Range("N132").Value = -14075.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1921.3125
$ws.Range("I20").Value = 1854.75
$ws.Range("J20").Value = 2121
$ws.Range("K20").Value = 1854.75
$ws.Range("L20").Value = 2121
$ws.Range("M20").Value = -1607.75
$ws.Range("N20").Value = -2615

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1767.6562
$ws.Range("I58").Value = 1284.4375
$ws.Range("J58").Value = 2250.875
$ws.Range("K58").Value = 1284.4375
$ws.Range("L58").Value = 2250.875
$ws.Range("M58").Value = -1081.4375
$ws.Range("N58").Value = -2656.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 4147.1816
$ws.Range("J94").Value = 4528.533
$ws.Range("L94").Value = 4528.533
$ws.Range("N94").Value = -5430.533

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2000
$ws.Range("I99").Value = 2000
$ws.Range("K99").Value = 2000
$ws.Range("M99").Value = -502

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3289.1428
$ws.Range("I122").Value = 3289.1428
$ws.Range("K122").Value = 9867.428400000001
$ws.Range("M122").Value = -7417.428400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2000
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1767.6562
$ws.Range("I136").Value = 1284.4375
$ws.Range("J136").Value = 2250.875
$ws.Range("K136").Value = 3853.3125
$ws.Range("L136").Value = 6752.625
$ws.Range("M136").Value = -1303.3125
$ws.Range("N136").Value = -11852.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 39533.332
$ws.Range("J140").Value = 39533.332
$ws.Range("L140").Value = 39533.332
$ws.Range("N140").Value = -49893.332

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 1379.2
$ws.Range("I118").Value = 701.5
$ws.Range("K118").Value = 2104.5
$ws.Range("M118").Value = -861.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2597.08
$ws.Range("I140").Value = 1401.35
$ws.Range("J140").Value = 7380
$ws.Range("K140").Value = 4204.049999999999
$ws.Range("L140").Value = 22140
$ws.Range("M140").Value = 975.9500000000007
$ws.Range("N140").Value = -32500

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7748.0757
$ws.Range("I70").Value = 10345.267
$ws.Range("J70").Value = 4360.4346
$ws.Range("K70").Value = 10345.267
$ws.Range("L70").Value = 4360.4346
$ws.Range("M70").Value = -10075.267
$ws.Range("N70").Value = -4900.4346

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 7748.0757
$ws.Range("I73").Value = 10345.267
$ws.Range("J73").Value = 4360.4346
$ws.Range("K73").Value = 10345.267
$ws.Range("L73").Value = 4360.4346
$ws.Range("M73").Value = -9409.267
$ws.Range("N73").Value = -6232.4346

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4447277.5
$ws.Range("I122").Value = 8335796.5
$ws.Range("J122").Value = 3255.7144
$ws.Range("K122").Value = 25007389.5
$ws.Range("L122").Value = 9767.143199999999
$ws.Range("M122").Value = -25004939.5
$ws.Range("N122").Value = -14667.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3212.9348
$ws.Range("I132").Value = 2223.8108
$ws.Range("J132").Value = 7279.3335
$ws.Range("K132").Value = 6671.432400000001
$ws.Range("L132").Value = 21838.0005
$ws.Range("M132").Value = -4141.432400000001
$ws.Range("N132").Value = -26898.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 58950
$ws.Range("J138").Value = 58950
$ws.Range("L138").Value = 58950
$ws.Range("N138").Value = -69230

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5094.2085
$ws.Range("I7").Value = 4950.9644
$ws.Range("J7").Value = 5294.75
$ws.Range("K7").Value = 4950.9644
$ws.Range("L7").Value = 5294.75
$ws.Range("M7").Value = -4838.9644
$ws.Range("N7").Value = -5518.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1036.4814
$ws.Range("I22").Value = 580
$ws.Range("J22").Value = 1140.2273
$ws.Range("K22").Value = 580
$ws.Range("L22").Value = 1140.2273
$ws.Range("M22").Value = -285
$ws.Range("N22").Value = -1730.2273

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1036.4814
$ws.Range("I27").Value = 580
$ws.Range("J27").Value = 1140.2273
$ws.Range("K27").Value = 580
$ws.Range("L27").Value = 1140.2273
$ws.Range("M27").Value = -473
$ws.Range("N27").Value = -1354.2273

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5304.9375
$ws.Range("I122").Value = 5392.1113
$ws.Range("J122").Value = 5192.857
$ws.Range("K122").Value = 16176.3339
$ws.Range("L122").Value = 15578.571
$ws.Range("M122").Value = -13726.3339
$ws.Range("N122").Value = -20478.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 5094.2085
$ws.Range("I126").Value = 4950.9644
$ws.Range("J126").Value = 5294.75
$ws.Range("K126").Value = 14852.8932
$ws.Range("L126").Value = 15884.25
$ws.Range("M126").Value = -12382.8932
$ws.Range("N126").Value = -20824.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7582101
$ws.Range("I132").Value = 3788.3403
$ws.Range("J132").Value = 26328454
$ws.Range("K132").Value = 11365.0209
$ws.Range("L132").Value = 78985362
$ws.Range("M132").Value = -8835.0209
$ws.Range("N132").Value = -78990422

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 61985.8
$ws.Range("J46").Value = 61985.8
$ws.Range("L46").Value = 61985.8
$ws.Range("N46").Value = -62447.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1548.0222
$ws.Range("I132").Value = 1416.65
$ws.Range("J132").Value = 2599
$ws.Range("K132").Value = 4249.950000000001
$ws.Range("L132").Value = 7797
$ws.Range("M132").Value = -1719.950000000001
$ws.Range("N132").Value = -12857

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H134").Value = 61985.8
$ws.Range("J134").Value = 61985.8
$ws.Range("L134").Value = 185957.4
$ws.Range("N134").Value = -191027.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 931.7273
$ws.Range("I136").Value = 720.02563
$ws.Range("J136").Value = 2583
$ws.Range("K136").Value = 2160.07689
$ws.Range("L136").Value = 7749
$ws.Range("M136").Value = 389.9231100000002
$ws.Range("N136").Value = -12849
